$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 76.08252666666667
$ws.Cells.Item(2, 8).Value = 228.24758
$ws.Cells.Item(2, 9).Value = 0.95878149807566
$ws.Cells.Item(2, 10).Value = 0.95878149807566
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 1.847798333333333
$ws.Cells.Item(2, 14).Value = 5.543395
$ws.Cells.Item(2, 15).Value = 0.05039680725746681
$ws.Cells.Item(2, 16).Value = 0.05039680725746681
$ws.Cells.Item(2, 17).Value = 140.5851659704556
$ws.Cells.Item(2, 18).Value = 1265.2664937341
$ws.Cells.Item(2, 19).Value = 0.04831952636054432
$ws.Cells.Item(2, 20).Value = 0.04831952636054432

$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 76.08252666666667
$ws.Cells.Item(3, 8).Value = 228.24758
$ws.Cells.Item(3, 9).Value = 0.95878149807566
$ws.Cells.Item(3, 10).Value = 0.95878149807566
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 0.4798556666666667
$ws.Cells.Item(3, 14).Value = 1.439567
$ws.Cells.Item(3, 15).Value = 0.0130875719001099
$ws.Cells.Item(3, 16).Value = 0.0130875719001099
$ws.Cells.Item(3, 17).Value = 36.50863155531778
$ws.Cells.Item(3, 18).Value = 328.57768399786
$ws.Cells.Item(3, 19).Value = 0.01254812179256028
$ws.Cells.Item(3, 20).Value = 0.01254812179256028

$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 76.08252666666667
$ws.Cells.Item(4, 8).Value = 228.24758
$ws.Cells.Item(4, 9).Value = 0.95878149807566
$ws.Cells.Item(4, 10).Value = 0.95878149807566
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 34.337334
$ws.Cells.Item(4, 14).Value = 103.012002
$ws.Cells.Item(4, 15).Value = 0.9365156208424232
$ws.Cells.Item(4, 16).Value = 0.9365156208424232
$ws.Cells.Item(4, 17).Value = 2612.47112971724
$ws.Cells.Item(4, 18).Value = 23512.24016745516
$ws.Cells.Item(4, 19).Value = 0.8979138499225554
$ws.Cells.Item(4, 20).Value = 0.8979138499225554

$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 0.9347409999999999
$ws.Cells.Item(5, 8).Value = 2.804223
$ws.Cells.Item(5, 9).Value = 0.01177947704364805
$ws.Cells.Item(5, 10).Value = 0.01177947704364805
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 1.847798333333333
$ws.Cells.Item(5, 14).Value = 5.543395
$ws.Cells.Item(5, 15).Value = 0.05039680725746681
$ws.Cells.Item(5, 16).Value = 0.05039680725746681
$ws.Cells.Item(5, 17).Value = 1.727212861898333
$ws.Cells.Item(5, 18).Value = 15.544915757085
$ws.Cells.Item(5, 19).Value = 0.0005936480341624857
$ws.Cells.Item(5, 20).Value = 0.0005936480341624857

$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 0.9347409999999999
$ws.Cells.Item(6, 8).Value = 2.804223
$ws.Cells.Item(6, 9).Value = 0.01177947704364805
$ws.Cells.Item(6, 10).Value = 0.01177947704364805
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 0.4798556666666667
$ws.Cells.Item(6, 14).Value = 1.439567
$ws.Cells.Item(6, 15).Value = 0.0130875719001099
$ws.Cells.Item(6, 16).Value = 0.0130875719001099
$ws.Cells.Item(6, 17).Value = 0.4485407657156666
$ws.Cells.Item(6, 18).Value = 4.036866891441
$ws.Cells.Item(6, 19).Value = 0.0001541647527544378
$ws.Cells.Item(6, 20).Value = 0.0001541647527544378

$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 0.9347409999999999
$ws.Cells.Item(7, 8).Value = 2.804223
$ws.Cells.Item(7, 9).Value = 0.01177947704364805
$ws.Cells.Item(7, 10).Value = 0.01177947704364805
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 34.337334
$ws.Cells.Item(7, 14).Value = 103.012002
$ws.Cells.Item(7, 15).Value = 0.9365156208424232
$ws.Cells.Item(7, 16).Value = 0.9365156208424232
$ws.Cells.Item(7, 17).Value = 32.096513920494
$ws.Cells.Item(7, 18).Value = 288.868625284446
$ws.Cells.Item(7, 19).Value = 0.01103166425673112
$ws.Cells.Item(7, 20).Value = 0.01103166425673112

$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 2.336085333333334
$ws.Cells.Item(8, 8).Value = 7.008256
$ws.Cells.Item(8, 9).Value = 0.02943902488069198
$ws.Cells.Item(8, 10).Value = 0.02943902488069198
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 1.847798333333333
$ws.Cells.Item(8, 14).Value = 5.543395
$ws.Cells.Item(8, 15).Value = 0.05039680725746681
$ws.Cells.Item(8, 16).Value = 0.05039680725746681
$ws.Cells.Item(8, 17).Value = 4.316614585457779
$ws.Cells.Item(8, 18).Value = 38.84953126912001
$ws.Cells.Item(8, 19).Value = 0.001483632862760003
$ws.Cells.Item(8, 20).Value = 0.001483632862760003

$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 2.336085333333334
$ws.Cells.Item(9, 8).Value = 7.008256
$ws.Cells.Item(9, 9).Value = 0.02943902488069198
$ws.Cells.Item(9, 10).Value = 0.02943902488069198
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 0.4798556666666667
$ws.Cells.Item(9, 14).Value = 1.439567
$ws.Cells.Item(9, 15).Value = 0.0130875719001099
$ws.Cells.Item(9, 16).Value = 0.0130875719001099
$ws.Cells.Item(9, 17).Value = 1.120983785016889
$ws.Cells.Item(9, 18).Value = 10.088854065152
$ws.Cells.Item(9, 19).Value = 0.0003852853547951805
$ws.Cells.Item(9, 20).Value = 0.0003852853547951806

$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 2.336085333333334
$ws.Cells.Item(10, 8).Value = 7.008256
$ws.Cells.Item(10, 9).Value = 0.02943902488069198
$ws.Cells.Item(10, 10).Value = 0.02943902488069198
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 34.337334
$ws.Cells.Item(10, 14).Value = 103.012002
$ws.Cells.Item(10, 15).Value = 0.9365156208424232
$ws.Cells.Item(10, 16).Value = 0.9365156208424232
$ws.Cells.Item(10, 17).Value = 80.214942343168
$ws.Cells.Item(10, 18).Value = 721.934481088512
$ws.Cells.Item(10, 19).Value = 0.02757010666313679
$ws.Cells.Item(10, 20).Value = 0.02757010666313679
